{"js": "// Replace the author name \"Adam Zabell\" with \"Randall Julian\".\n// Do the two word-level replacements separately (instead of replacing the\n// whole phrase) so that the existing run layout / formatting of the\n// paragraph (e.g. the separate run holding the space) is preserved as\n// closely as possible, matching how the original document splits the\n// author's first and last name into distinct runs.\n\nconst firstNameResults = context.document.body.search(\"Adam\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nfirstNameResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < firstNameResults.items.length; i++) {\n  firstNameResults.items[i].insertText(\"Randall\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\nconst lastNameResults = context.document.body.search(\"Zabell\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nlastNameResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < lastNameResults.items.length; i++) {\n  lastNameResults.items[i].insertText(\"Julian\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the author name \"Adam Zabell\" with \"Randall Julian\".\n# Use Find & Replace (the standard Word automation idiom) against the\n# whole document body, doing the first-name and last-name swaps as two\n# separate whole-word, case-sensitive replacements so the unrelated\n# space between them is left untouched.\n\n$d = $word.ActiveDocument\n\n$firstName = $d.Content.Find\n$firstName.Execute(\n    \"Adam\",    # FindText\n    $true,     # MatchCase\n    $true,     # MatchWholeWord\n    $false,    # MatchWildcards\n    $false,    # MatchSoundsLike\n    $false,    # MatchAllWordForms\n    $true,     # Forward\n    1,         # Wrap (wdFindContinue)\n    $false,    # Format\n    \"Randall\", # ReplaceWith\n    2          # Replace (wdReplaceAll)\n) | Out-Null\n\n$lastName = $d.Content.Find\n$lastName.Execute(\n    \"Zabell\",  # FindText\n    $true,     # MatchCase\n    $true,     # MatchWholeWord\n    $false,    # MatchWildcards\n    $false,    # MatchSoundsLike\n    $false,    # MatchAllWordForms\n    $true,     # Forward\n    1,         # Wrap (wdFindContinue)\n    $false,    # Format\n    \"Julian\",  # ReplaceWith\n    2          # Replace (wdReplaceAll)\n) | Out-Null\n"}
